$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.727054357528687
$ws.Range("B1").Value = 4.954568862915039
$ws.Range("C1").Value = 4.684307098388672
$ws.Range("D1").Value = 7.788054466247559
$ws.Range("E1").Value = 5.472336769104004
